$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77 - this shifts existing rows 77..118 down to 78..119
# (matches the diff: a new weekly record is inserted, everything after shifts down one row)
$ws.Rows.Item(77).Insert()

# Populate the new row 77 with the inserted record's data
$ws.Range("A77").Value = 8
$ws.Range("B77").Value = "Terminal La Palmera de La Serena"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value = 45016
$ws.Range("D77").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E77").Value = 4
$ws.Range("F77").Value = 100114007
$ws.Range("G77").Value = "Jengibre"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 430
$ws.Range("K77").Value = 16000
$ws.Range("L77").Value = 17000
$ws.Range("M77").Value = 16500
$ws.Range("N77").Value = "$/caja 13 kilos"
$ws.Range("O77").Value = "Perú"
$ws.Range("P77").Value = 1269
$ws.Range("Q77").Value = 13
$ws.Range("R77").Value = "Hortaliza"
